$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 42/43: swap Coin name and Link (RenderToken <-> NEARProtocol)
# Row 2
$ws.Range("D2").Value = "43.295.76"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").Value = "2.550.75"
$ws.Range("E3").Value = "  +0.69%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.61"
$ws.Range("E5").Value = "  +4.30%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.52"
$ws.Range("E6").Value = "  -2.16%  "

# Row 7
$ws.Range("E7").Value = "  -0.42%  "

# Row 8
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -1.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.75"
$ws.Range("E10").Value = "  -0.77%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("E11").Value = "  +0.30%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.68"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13
$ws.Range("E13").Value = "  +0.52%  "

# Row 14
$ws.Range("D14").Value = "2.942.97"
$ws.Range("E14").Value = "  +0.72%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.65"
$ws.Range("E15").Value = "  +3.58%  "

# Row 16
$ws.Range("D16").Value = "2.524.12"
$ws.Range("E16").Value = "  -0.77%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.857"
$ws.Range("E17").Value = "  -1.22%  "

# Row 18
$ws.Range("D18").Value = "43.194.57"
$ws.Range("E18").Value = "  +0.78%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("E19").Value = "  +1.57%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.71"
$ws.Range("E20").Value = "  +3.55%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0974"
$ws.Range("E21").Value = "  -0.81%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.67"
$ws.Range("E22").Value = "  -1.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.21"
$ws.Range("E23").Value = "  +0.05%  "

# Row 24
$ws.Range("E24").Value = "  +1.68%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.04"
$ws.Range("E25").Value = "  -0.36%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.20"
$ws.Range("E26").Value = "  +1.24%  "

# Row 27
$ws.Range("E27").Value = "  -0.12%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.43"
$ws.Range("E28").Value = "  +3.55%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.10"
$ws.Range("E29").Value = "  +4.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.28"
$ws.Range("E30").Value = "  -1.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.13"
$ws.Range("E31").Value = "  +0.67%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.73"
$ws.Range("E32").Value = "  -1.34%  "

# Row 33
$ws.Range("E33").Value = "  +2.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.37"
$ws.Range("E34").Value = "  +1.36%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.22"
$ws.Range("E35").Value = "  +5.40%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0795"
$ws.Range("E36").Value = "  +0.18%  "

# Row 37
$ws.Range("E37").Value = "  +0.14%  "

# Row 38
$ws.Range("E38").Value = "  -2.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "24.75"
$ws.Range("E39").Value = "  +1.99%  "

# Row 40
$ws.Range("E40").Value = "  -0.78%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.28"
$ws.Range("E41").Value = "  +9.44%  "

# Row 42
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.40"
$ws.Range("E42").Value = "  -1.22%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.85"
$ws.Range("E43").Value = "  -1.39%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0306"
$ws.Range("E44").Value = "  +0.77%  "

# Row 45
$ws.Range("E45").Value = "  +0.22%  "

# Row 46
$ws.Range("D46").Value = "2.024.06"
$ws.Range("E46").Value = "  -1.11%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.32"
$ws.Range("E47").Value = "  +0.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.90"
$ws.Range("E48").Value = "  -0.65%  "

# Row 49
$ws.Range("D49").Value = "2.795.70"
$ws.Range("E49").Value = "  +0.58%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.95"
$ws.Range("E50").Value = "  +2.94%  "

# Row 51
$ws.Range("E51").Value = "  +0.34%  "
